$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ((Intercept))
$ws.Range("B2").Value = 1396.011396
$ws.Range("D2").Value = 2.972005
$ws.Range("E2").Value = 0.086103

# Row 3 (household_group_collapsed)
$ws.Range("B3").Value = 818.935035
$ws.Range("D3").Value = 0.871726
$ws.Range("E3").Value = 0.419649

# Row 4 (Residuals)
$ws.Range("B4").Value = 104747.660202
$ws.Range("C4").Value = 223

# Row 5 (SM-Control)
$ws.Range("G5").Value = 5.179519
$ws.Range("H5").Value = -4.285879
$ws.Range("I5").Value = 14.644916
$ws.Range("J5").Value = 0.401663

# Row 6 (SM + Traps-Control)
$ws.Range("G6").Value = 4.735163
$ws.Range("H6").Value = -5.457066
$ws.Range("I6").Value = 14.927392
$ws.Range("J6").Value = 0.517438

# Row 7 (SM + Traps-SM)
$ws.Range("G7").Value = -0.444356
$ws.Range("H7").Value = -8.149787
$ws.Range("I7").Value = 7.261075
$ws.Range("J7").Value = 0.989846
